$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 51.333332
$ws.Range("I11").Value = 51.333332
$ws.Range("K11").Value = 51.333332
$ws.Range("M11").Value = 88.666668

$ws.Range("H19").Value = 1527.76
$ws.Range("I19").Value = 1314.3636
$ws.Range("J19").Value = 1695.4286
$ws.Range("K19").Value = 1314.3636
$ws.Range("L19").Value = 1695.4286
$ws.Range("M19").Value = -1139.3636
$ws.Range("N19").Value = -2045.4286

$ws.Range("H33").Value = 92.82353000000001
$ws.Range("I33").Value = 56.625
$ws.Range("K33").Value = 56.625
$ws.Range("M33").Value = 172.375

$ws.Range("H112").Value = 2799.6667
$ws.Range("J112").Value = 2799.6667
$ws.Range("L112").Value = 8399.000100000001
$ws.Range("N112").Value = -10615.0001

$ws.Range("H137").Value = 36455.965
$ws.Range("I137").Value = 1002
$ws.Range("J137").Value = 37722.18
$ws.Range("K137").Value = 3006
$ws.Range("L137").Value = 113166.54
$ws.Range("M137").Value = -456
$ws.Range("N137").Value = -118266.54

$ws.Range("H138").Value = 2425.1135
$ws.Range("I138").Value = 2645.125
$ws.Range("J138").Value = 2299.3928
$ws.Range("K138").Value = 7935.375
$ws.Range("L138").Value = 6898.178400000001
$ws.Range("M138").Value = -2795.375
$ws.Range("N138").Value = -17178.1784

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 1079228.2
$ws.Range("I141").Value = 1475523
$ws.Range("J141").Value = 3571.1428
$ws.Range("K141").Value = 4426569
$ws.Range("L141").Value = 10713.4284
$ws.Range("M141").Value = -4421389
$ws.Range("N141").Value = -21073.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2225.3098
$ws.Range("I32").Value = 1746.8392
$ws.Range("J32").Value = 4011.6
$ws.Range("K32").Value = 1746.8392
$ws.Range("L32").Value = 4011.6
$ws.Range("M32").Value = -1459.8392
$ws.Range("N32").Value = -4585.6

$ws.Range("H45").Value = 2936.8845
$ws.Range("I45").Value = 3452.1428
$ws.Range("J45").Value = 2747.0527
$ws.Range("K45").Value = 3452.1428
$ws.Range("L45").Value = 2747.0527
$ws.Range("M45").Value = -3075.1428
$ws.Range("N45").Value = -3501.0527

$ws.Range("H61").Value = 3603.762
$ws.Range("I61").Value = 2688.1
$ws.Range("J61").Value = 4436.1816
$ws.Range("K61").Value = 2688.1
$ws.Range("L61").Value = 4436.1816
$ws.Range("M61").Value = -2476.1
$ws.Range("N61").Value = -4860.1816

$ws.Range("H74").Value = 1880.9
$ws.Range("J74").Value = 4633.3335
$ws.Range("L74").Value = 4633.3335
$ws.Range("N74").Value = -6381.3335

$ws.Range("H77").Value = 1880.9
$ws.Range("J77").Value = 4633.3335
$ws.Range("L77").Value = 23166.6675
$ws.Range("N77").Value = -31902.6675

$ws.Range("H122").Value = 49158.75
$ws.Range("I122").Value = 77674
$ws.Range("J122").Value = 1633.3334
$ws.Range("K122").Value = 233022
$ws.Range("L122").Value = 4900.0002
$ws.Range("M122").Value = -230572
$ws.Range("N122").Value = -9800.0002

$ws.Range("H132").Value = 1942.6316
$ws.Range("I132").Value = 1740.38
$ws.Range("J132").Value = 3387.2856
$ws.Range("K132").Value = 5221.14
$ws.Range("L132").Value = 10161.8568
$ws.Range("M132").Value = -2691.14
$ws.Range("N132").Value = -15221.8568

$ws.Range("H136").Value = 3603.762
$ws.Range("I136").Value = 2688.1
$ws.Range("J136").Value = 4436.1816
$ws.Range("K136").Value = 8064.299999999999
$ws.Range("L136").Value = 13308.5448
$ws.Range("M136").Value = -5514.299999999999
$ws.Range("N136").Value = -18408.5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2183.4517
$ws.Range("I105").Value = 2189.2068
$ws.Range("K105").Value = 2189.2068
$ws.Range("M105").Value = -442.2067999999999

$ws.Range("H107").Value = 5262.3335
$ws.Range("I107").Value = 4993
$ws.Range("K107").Value = 4993
$ws.Range("M107").Value = -3073

$ws.Range("H134").Value = 12668.8
$ws.Range("I134").Value = 16612.572
$ws.Range("K134").Value = 49837.716
$ws.Range("M134").Value = -47302.716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1308.4474
$ws.Range("I31").Value = 804.9091
$ws.Range("J31").Value = 2000.8125
$ws.Range("K31").Value = 804.9091
$ws.Range("L31").Value = 2000.8125
$ws.Range("M31").Value = -509.9091
$ws.Range("N31").Value = -2590.8125

$ws.Range("H34").Value = 1308.4474
$ws.Range("I34").Value = 804.9091
$ws.Range("J34").Value = 2000.8125
$ws.Range("K34").Value = 804.9091
$ws.Range("L34").Value = 2000.8125
$ws.Range("M34").Value = -602.9091
$ws.Range("N34").Value = -2404.8125

$ws.Range("H99").Value = 2920
$ws.Range("I99").Value = 2775
$ws.Range("K99").Value = 2775
$ws.Range("M99").Value = -1277

$ws.Range("H126").Value = 2920
$ws.Range("I126").Value = 2775
$ws.Range("K126").Value = 8325
$ws.Range("M126").Value = -5855

$ws.Range("H132").Value = 2977.3914
$ws.Range("I132").Value = 2637
$ws.Range("J132").Value = 3419.9
$ws.Range("K132").Value = 7911
$ws.Range("L132").Value = 10259.7
$ws.Range("M132").Value = -5381
$ws.Range("N132").Value = -15319.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 45.235294
$ws.Range("I12").Value = 16.636364
$ws.Range("K12").Value = 49.909092
$ws.Range("M12").Value = 123.090908

$ws.Range("H17").Value = 250000800
$ws.Range("I17").Value = 333334080
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 1000002240
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -1000002071
$ws.Range("N17").Value = -3338

$ws.Range("H131").Value = 6340339.5
$ws.Range("I131").Value = 250000510
$ws.Range("J131").Value = 11503.571
$ws.Range("K131").Value = 750001530
$ws.Range("L131").Value = 34510.713
$ws.Range("M131").Value = -749996490
$ws.Range("N131").Value = -44590.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 196.46666
$ws.Range("I2").Value = 212
$ws.Range("K2").Value = 212
$ws.Range("M2").Value = -99

$ws.Range("H102").Value = 4891.7144
$ws.Range("I102").Value = 4891.7144
$ws.Range("K102").Value = 4891.7144
$ws.Range("M102").Value = -3269.7144

$ws.Range("H122").Value = 3905.5
$ws.Range("J122").Value = 4116.4443
$ws.Range("L122").Value = 12349.3329
$ws.Range("N122").Value = -17249.3329

$ws.Range("H126").Value = 1770369.9
$ws.Range("I126").Value = 2927175.2
$ws.Range("J126").Value = 79654.38
$ws.Range("K126").Value = 8781525.600000001
$ws.Range("L126").Value = 238963.14
$ws.Range("M126").Value = -8779055.600000001
$ws.Range("N126").Value = -243903.14

$ws.Range("H132").Value = 1014669.8
$ws.Range("I132").Value = 1540120.2
$ws.Range("J132").Value = 4188.385
$ws.Range("K132").Value = 4620360.6
$ws.Range("L132").Value = 12565.155
$ws.Range("M132").Value = -4617830.6
$ws.Range("N132").Value = -17625.155

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 272.69232
$ws.Range("I55").Value = 231.36363
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 231.36363
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -58.36363
$ws.Range("N55").Value = -846

$ws.Range("H132").Value = 2434.4348
$ws.Range("I132").Value = 970.5
$ws.Range("J132").Value = 4031.4546
$ws.Range("K132").Value = 2911.5
$ws.Range("L132").Value = 12094.3638
$ws.Range("M132").Value = -381.5
$ws.Range("N132").Value = -17154.3638

$ws.Range("H136").Value = 3840.5293
$ws.Range("I136").Value = 2662.7273
$ws.Range("J136").Value = 5999.8335
$ws.Range("K136").Value = 7988.1819
$ws.Range("L136").Value = 17999.5005
$ws.Range("M136").Value = -5438.1819
$ws.Range("N136").Value = -23099.5005

$ws.Range("H139").Value = 30900
$ws.Range("J139").Value = 30900
$ws.Range("L139").Value = 30900
$ws.Range("N139").Value = -41180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5818.1816
$ws.Range("J96").Value = 5818.1816
$ws.Range("L96").Value = 5818.1816
$ws.Range("N96").Value = -8564.1816

$ws.Range("H107").Value = 1566.9
$ws.Range("I107").Value = 1650
$ws.Range("J107").Value = 1546.125
$ws.Range("K107").Value = 4950
$ws.Range("L107").Value = 4638.375
$ws.Range("M107").Value = -3030
$ws.Range("N107").Value = -8478.375

$ws.Range("H126").Value = 4276.9614
$ws.Range("I126").Value = 3570.0527
$ws.Range("J126").Value = 6195.7144
$ws.Range("K126").Value = 10710.1581
$ws.Range("L126").Value = 18587.1432
$ws.Range("M126").Value = -8240.158100000001
$ws.Range("N126").Value = -23527.1432

$ws.Range("H132").Value = 2827.6667
$ws.Range("I132").Value = 2441.625
$ws.Range("J132").Value = 3599.75
$ws.Range("K132").Value = 7324.875
$ws.Range("L132").Value = 10799.25
$ws.Range("M132").Value = -4794.875
$ws.Range("N132").Value = -15859.25

$ws.Range("H136").Value = 25256212
$ws.Range("I136").Value = 55558868
$ws.Range("K136").Value = 166676604
$ws.Range("M136").Value = -166674054

Write-Host "Applied 42 row updates across 8 sheets"
